$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = "Dr. Hend Mahmoud, Dr. Majorelle Magdy, Dr. Servinaz Sayed Mohammad, Dr. Eman Tantawi"
$ws.Range("G3").Value = "Dr. Amira Sobhy, Dr. Menna tuâ€™Allah Medhat, Dr. Veronia Rafat, Dr. Asmaa Reda, Dr. Eman Tantawi"
$ws.Range("G4").Value = "Dr. Amira Sobhy, Dr. Rana Abo-Zaid, Dr. Hend Mahmoud, Dr. Veronia Rafat, Dr. Asmaa Reda, Dr. Eman Tantawi"
$ws.Range("G5").Value = "Dr. Nesma, Dr. Servinaz Sayed Mohammad, Dr. Eman Tantawi, Dr. Nourhan Mahmoud, Dr. Hanan Ragab, Dr. Veronia Rafat, Dr. Mohammad El-Tanany, Dr. Hend Mahmoud"
$ws.Range("G6").Value = "Dr. Amira Sobhy, Dr. Servinaz Sayed Mohammad, Dr. Hend Mahmoud, Dr. Menna tuâ€™Allah Medhat, Dr. Nourhan Mahmoud, Dr. Nahla Nagiub, Dr. Veronia Rafat, Dr. Gehan Adel, Dr. Asmaa Reda, Dr. Eman Tantawi"
$ws.Range("G7").Value = "Dr. Amira Sobhy, Dr. Rana Abo-Zaid, Dr. Servinaz Sayed Mohammad, Dr. Menna tuâ€™Allah Medhat, Dr. Eman Tantawi, Dr. Veronia Rafat, Dr. Gehan Adel, Dr. Asmaa Reda, Dr. Hend Mahmoud"
$ws.Range("G8").Value = "Dr. Shimaa Ahmad Mekki, Administrator, Dr. Manar Montaser, Dr. Majorelle Magdy, Dr. Asmaa Reda, Dr. Eman Tantawi"
$ws.Range("G9").Value = "Dr. Amira Sobhy, Dr. Rana Abo-Zaid, Dr. Manar Montaser, Dr. Majorelle Magdy, Dr. Menna tuâ€™Allah Medhat, Dr. Gehan Adel, Dr. Asmaa Reda, Dr. Hend Mahmoud"
$ws.Range("G10").Value = "Dr. Servinaz Sayed Mohammad, Dr. Shimaa Ahmad Mekki, Dr. Rana Abo-Zaid, Dr. Sara Wael, Dr. Heba Mahmoud Ali, Dr. Alshimaa Atef, Dr. Gehan Adel"
$ws.Range("G12").Value = "Dr. Salma El-Gendy, Administrator"
$ws.Range("G13").Value = "Dr. Mariam Nour El-Din, D Wessam Atef, Dr. Shimaa Ashraf, Dr. Omnia Mohammad, Dr. Safa Hany"
$ws.Range("G14").Value = "Dr. Shimaa Ashraf, Dr. Safa Hany"
$ws.Range("G15").Value = "Dr. Amal Awwad, D Wessam Atef"
$ws.Range("G17").Value = "Dr. Eman M. Abo-Sakaya, Dr. Basma Hamed, Dr. Esraa Mostafa, Dr. Nourhan Osama, Dr. Dina Adel, Dr. Yasmeena Fattoh, Dr. Arwa Al-Sayed, Dr. Sarah Abdelmohsen, Dr. Madeha Saeed, Dr. Marwa Mustafa"
$ws.Range("G19").Value = "D Mariam E. Mohammad, Dr. Sarah Mahdy"
$ws.Range("G24").Value = "Dr. Aya Emad, Dr. Maryam Ashraf, Dr. Neveen Nashaat, Dr. Youstina Magdy, Dr. Yasmin, Dr. Marina Atef, Dr. Wafaa Ebida, Dr. Ola Abd Al-Fattah, Dr. Salma Hassan, Dr. Remon, Dr. Monica"
$ws.Range("G25").Value = "Dr. Aya Emad, Dr. Youstina Magdy, Dr. Marina Atef, Dr. Ola Abd Al-Fattah, Dr. Eman Samir Gabry, Dr. Remon, Dr. Abdullah El-Agrody"
$ws.Range("G27").Value = "Dr. Neveen Nashaat, Dr. Eman Mohammad Al, Dr. Yasmin, Dr. Wafaa Ebida, Dr. Ola Abd Al-Fattah, Dr. Salma Hassan, Dr. Eman Samir Gabry, Dr. Remon"
$ws.Range("G28").Value = "Dr. Aya Hanafy, Dr. Nardine, Dr. Neveen Nashaat, Dr. Wafaa Ebida, Dr. Salma Hassan, Dr. Eman Samir Gabry, Dr. Remon, Dr. Abdullah El-Agrody"
$ws.Range("G29").Value = "Dr. Naema Gomaa, Dr. Neveen Nashaat, Dr. Ola Abd Al-Fattah, Dr. Eman Samir Gabry, Dr. Remon, Dr. Monica"
$ws.Range("G30").Value = "Dr. Amira Sobhy, Dr. Rana Abo-Zaid, Dr. Hend Mahmoud, Dr. Veronia Rafat, Dr. Asmaa Reda, Dr. Eman Tantawi"
$ws.Range("G31").Value = "Dr. Amira Sobhy, Dr. Menna tuâ€™Allah Medhat, Dr. Veronia Rafat, Dr. Asmaa Reda, Dr. Eman Tantawi"
$ws.Range("G32").Value = "Dr. Amira Sobhy, Dr. Rana Abo-Zaid, Dr. Hend Mahmoud, Dr. Veronia Rafat, Dr. Asmaa Reda, Dr. Eman Tantawi"
$ws.Range("G33").Value = "Dr. Nesma, Dr. Servinaz Sayed Mohammad, Dr. Eman Tantawi, Dr. Nourhan Mahmoud, Dr. Hanan Ragab, Dr. Veronia Rafat, Dr. Mohammad El-Tanany, Dr. Hend Mahmoud"
$ws.Range("G34").Value = "Dr. Amira Sobhy, Dr. Servinaz Sayed Mohammad, Dr. Hend Mahmoud, Dr. Menna tuâ€™Allah Medhat, Dr. Nourhan Mahmoud, Dr. Nahla Nagiub, Dr. Veronia Rafat, Dr. Gehan Adel, Dr. Asmaa Reda, Dr. Eman Tantawi"
$ws.Range("G35").Value = "Dr. Amira Sobhy, Dr. Rana Abo-Zaid, Dr. Servinaz Sayed Mohammad, Dr. Menna tuâ€™Allah Medhat, Dr. Eman Tantawi, Dr. Veronia Rafat, Dr. Gehan Adel, Dr. Asmaa Reda, Dr. Hend Mahmoud"
$ws.Range("G36").Value = "Dr. Shimaa Ahmad Mekki, Administrator, Dr. Manar Montaser, Dr. Majorelle Magdy, Dr. Asmaa Reda, Dr. Eman Tantawi"
$ws.Range("G37").Value = "Dr. Amira Sobhy, Dr. Rana Abo-Zaid, Dr. Manar Montaser, Dr. Majorelle Magdy, Dr. Menna tuâ€™Allah Medhat, Dr. Gehan Adel, Dr. Asmaa Reda, Dr. Hend Mahmoud"
$ws.Range("G38").Value = "Dr. Servinaz Sayed Mohammad, Dr. Shimaa Ahmad Mekki, Dr. Rana Abo-Zaid, Dr. Sara Wael, Dr. Heba Mahmoud Ali, Dr. Alshimaa Atef, Dr. Gehan Adel"
$ws.Range("G40").Value = "Dr. Salma El-Gendy, Administrator"
$ws.Range("G41").Value = "Dr. Mariam Nour El-Din, D Wessam Atef, Dr. Shimaa Ashraf, Dr. Omnia Mohammad, Dr. Safa Hany"
$ws.Range("G42").Value = "Dr. Shimaa Ashraf, Dr. Safa Hany"
$ws.Range("G43").Value = "Dr. Amal Awwad, D Wessam Atef"
$ws.Range("G45").Value = "Dr. Eman M. Abo-Sakaya, Dr. Basma Hamed, Dr. Esraa Mostafa, Dr. Nourhan Osama, Dr. Dina Adel, Dr. Yasmeena Fattoh, Dr. Arwa Al-Sayed, Dr. Sarah Abdelmohsen, Dr. Madeha Saeed, Dr. Marwa Mustafa"
$ws.Range("G47").Value = "D Mariam E. Mohammad, Dr. Sarah Mahdy"
$ws.Range("G52").Value = "Dr. Aya Emad, Dr. Maryam Ashraf, Dr. Neveen Nashaat, Dr. Youstina Magdy, Dr. Yasmin, Dr. Marina Atef, Dr. Wafaa Ebida, Dr. Ola Abd Al-Fattah, Dr. Salma Hassan, Dr. Remon, Dr. Monica"
$ws.Range("G53").Value = "Dr. Aya Emad, Dr. Youstina Magdy, Dr. Marina Atef, Dr. Ola Abd Al-Fattah, Dr. Eman Samir Gabry, Dr. Remon, Dr. Abdullah El-Agrody"
$ws.Range("G55").Value = "Dr. Neveen Nashaat, Dr. Eman Mohammad Al, Dr. Yasmin, Dr. Wafaa Ebida, Dr. Ola Abd Al-Fattah, Dr. Salma Hassan, Dr. Eman Samir Gabry, Dr. Remon"
$ws.Range("G56").Value = "Dr. Aya Hanafy, Dr. Nardine, Dr. Neveen Nashaat, Dr. Wafaa Ebida, Dr. Salma Hassan, Dr. Eman Samir Gabry, Dr. Remon, Dr. Abdullah El-Agrody"
$ws.Range("G57").Value = "Dr. Naema Gomaa, Dr. Neveen Nashaat, Dr. Ola Abd Al-Fattah, Dr. Eman Samir Gabry, Dr. Remon, Dr. Monica"
